$wb = $excel.ActiveWorkbook

# ----- Sheet: ALC -----
$ws = $wb.Worksheets.Item("ALC")

# Row 86
$ws.Range("H86").Value = 4499.5
$ws.Range("I86").Value = 4499.5
$ws.Range("K86").Value = 4499.5
$ws.Range("M86").Value = -3376.5
# Row 89
$ws.Range("H89").Value = 4499.5
$ws.Range("I89").Value = 4499.5
$ws.Range("K89").Value = 22497.5
$ws.Range("M89").Value = -16881.5
# Row 112
$ws.Range("H112").Value = 1446.409
$ws.Range("J112").Value = 1515.238
$ws.Range("L112").Value = 4545.714
$ws.Range("N112").Value = -6761.714
# Row 113
$ws.Range("H113").Value = 6999.5
$ws.Range("I113").Value = 6999
$ws.Range("K113").Value = 6999
$ws.Range("M113").Value = -3745
# Row 116
$ws.Range("H116").Value = 3716.2666
$ws.Range("I116").Value = 2962
$ws.Range("J116").Value = 4376.25
$ws.Range("K116").Value = 2962
$ws.Range("L116").Value = 4376.25
$ws.Range("M116").Value = 480
$ws.Range("N116").Value = -11260.25
# Row 121
$ws.Range("H121").Value = 1770.2727
$ws.Range("J121").Value = 1937.8
$ws.Range("L121").Value = 5813.4
$ws.Range("N121").Value = -9307.4
# Row 138
$ws.Range("H138").Value = 4243.3
$ws.Range("J138").Value = 4625
$ws.Range("L138").Value = 13875
$ws.Range("N138").Value = -24155

# ----- Sheet: ARM -----
$ws = $wb.Worksheets.Item("ARM")

# Row 74
$ws.Range("H74").Value = 9998080
$ws.Range("I74").Value = 13329107
$ws.Range("K74").Value = 13329107
$ws.Range("M74").Value = -13328233
# Row 77
$ws.Range("H77").Value = 9998080
$ws.Range("I77").Value = 13329107
$ws.Range("K77").Value = 66645535
$ws.Range("M77").Value = -66641167
# Row 97
$ws.Range("H97").Value = 1782.4
$ws.Range("J97").Value = 2944
$ws.Range("L97").Value = 2944
$ws.Range("N97").Value = -3936
# Row 133
$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").Value = ""
# Row 135
$ws.Range("H135").Value = 100000
$ws.Range("J135").Value = 100000
$ws.Range("L135").Value = 100000
$ws.Range("N135").Value = -110140

# ----- Sheet: BSM -----
$ws = $wb.Worksheets.Item("BSM")

# Row 20
$ws.Range("H20").Value = 2496
$ws.Range("I20").Value = 2496
$ws.Range("K20").Value = 2496
$ws.Range("M20").Value = -2249
# Row 24
$ws.Range("H24").Value = 10638.667
$ws.Range("I24").Value = 1916
$ws.Range("J24").Value = 15000
$ws.Range("K24").Value = 1916
$ws.Range("L24").Value = 15000
$ws.Range("M24").Value = -1681
$ws.Range("N24").Value = -15470
# Row 33
$ws.Range("H33").Value = 10021
$ws.Range("I33").Value = 10021
$ws.Range("K33").Value = 10021
$ws.Range("M33").Value = -9685
# Row 86
$ws.Range("H86").Value = 3736.1428
$ws.Range("I86").Value = 2886.4614
$ws.Range("J86").Value = 5116.875
$ws.Range("K86").Value = 2886.4614
$ws.Range("L86").Value = 5116.875
$ws.Range("M86").Value = -1763.4614
$ws.Range("N86").Value = -7362.875
# Row 89
$ws.Range("H89").Value = 3736.1428
$ws.Range("I89").Value = 2886.4614
$ws.Range("J89").Value = 5116.875
$ws.Range("K89").Value = 14432.307
$ws.Range("L89").Value = 25584.375
$ws.Range("M89").Value = -8816.307000000001
$ws.Range("N89").Value = -36816.375
# Row 105
$ws.Range("H105").Value = 2498.6667
$ws.Range("I105").Value = 2362.0908
$ws.Range("K105").Value = 2362.0908
$ws.Range("M105").Value = -615.0907999999999
# Row 132
$ws.Range("H132").Value = 124400
$ws.Range("J132").Value = 124400
$ws.Range("L132").Value = 124400
$ws.Range("N132").Value = -134520
# Row 134
$ws.Range("H134").Value = 1543.5555
$ws.Range("I134").Value = 1442.3125
$ws.Range("K134").Value = 4326.9375
$ws.Range("M134").Value = -1791.9375
# Row 135
$ws.Range("H135").Value = 45884.332
$ws.Range("J135").Value = 45884.332
$ws.Range("L135").Value = 45884.332
$ws.Range("N135").Value = -56024.332
# Row 137
$ws.Range("H137").Value = 99900
$ws.Range("J137").Value = 99900
$ws.Range("L137").Value = 99900
$ws.Range("N137").Value = -110100
# Row 138
$ws.Range("H138").Value = 122566
$ws.Range("J138").Value = 122566
$ws.Range("L138").Value = 122566
$ws.Range("N138").Value = -132846

# ----- Sheet: CUL -----
$ws = $wb.Worksheets.Item("CUL")

# Row 38
$ws.Range("H38").Value = 147
$ws.Range("I38").Value = 50
$ws.Range("J38").Value = 171.25
$ws.Range("K38").Value = 150
$ws.Range("L38").Value = 513.75
$ws.Range("M38").Value = 197
$ws.Range("N38").Value = -1207.75
# Row 75
$ws.Range("J75").Value = 1310
$ws.Range("L75").Value = 3930
$ws.Range("N75").Value = -5926
# Row 78
$ws.Range("J78").Value = 1310
$ws.Range("L78").Value = 11790
$ws.Range("N78").Value = -21774
# Row 117
$ws.Range("H117").Value = 1387.6666
$ws.Range("J117").Value = 1531.5
$ws.Range("L117").Value = 4594.5
$ws.Range("N117").Value = -11478.5
# Row 129
$ws.Range("H129").Value = 1977
$ws.Range("I129").Value = 1949
$ws.Range("J129").Value = 2033
$ws.Range("K129").Value = 5847
$ws.Range("L129").Value = 6099
$ws.Range("M129").Value = -847
$ws.Range("N129").Value = -16099
# Row 138
$ws.Range("H138").Value = 5033.8335
$ws.Range("I138").Value = 4240.6
$ws.Range("J138").Value = 9000
$ws.Range("K138").Value = 12721.8
$ws.Range("L138").Value = 27000
$ws.Range("M138").Value = -7581.800000000001
$ws.Range("N138").Value = -37280

# ----- Sheet: GSM -----
$ws = $wb.Worksheets.Item("GSM")

# Row 80
$ws.Range("H80").Value = 3565.3333
$ws.Range("I80").Value = 2763.8333
$ws.Range("J80").Value = 5168.3335
$ws.Range("K80").Value = 2763.8333
$ws.Range("L80").Value = 5168.3335
$ws.Range("M80").Value = -1765.8333
$ws.Range("N80").Value = -7164.3335
# Row 83
$ws.Range("H83").Value = 3565.3333
$ws.Range("I83").Value = 2763.8333
$ws.Range("J83").Value = 5168.3335
$ws.Range("K83").Value = 13819.1665
$ws.Range("L83").Value = 25841.6675
$ws.Range("M83").Value = -8827.166499999999
$ws.Range("N83").Value = -35825.6675
# Row 126
$ws.Range("H126").Value = 3872.625
$ws.Range("I126").Value = 3799.5
$ws.Range("J126").Value = 3897
$ws.Range("K126").Value = 11398.5
$ws.Range("L126").Value = 11691
$ws.Range("M126").Value = -8928.5
$ws.Range("N126").Value = -16631
# Row 128
$ws.Range("H128").Value = 0
$ws.Range("J128").Value = 0
$ws.Range("L128").Value = 0
$ws.Range("N128").Value = ""

# ----- Sheet: LTW -----
$ws = $wb.Worksheets.Item("LTW")

# Row 22
$ws.Range("H22").Value = 930.5
$ws.Range("I22").Value = 900
$ws.Range("J22").Value = 961
$ws.Range("K22").Value = 900
$ws.Range("L22").Value = 961
$ws.Range("M22").Value = -605
$ws.Range("N22").Value = -1551
# Row 27
$ws.Range("H27").Value = 930.5
$ws.Range("I27").Value = 900
$ws.Range("J27").Value = 961
$ws.Range("K27").Value = 900
$ws.Range("L27").Value = 961
$ws.Range("M27").Value = -793
$ws.Range("N27").Value = -1175
# Row 46
$ws.Range("H46").Value = 924.5
$ws.Range("I46").Value = 899.6667
$ws.Range("K46").Value = 899.6667
$ws.Range("M46").Value = -711.6667
# Row 93
$ws.Range("H93").Value = 2915
$ws.Range("J93").Value = 2887
$ws.Range("L93").Value = 2887
$ws.Range("N93").Value = -5383
# Row 99
$ws.Range("H99").Value = 0
$ws.Range("I99").Value = 0
$ws.Range("K99").Value = 0
$ws.Range("M99").Value = ""
# Row 132
$ws.Range("H132").Value = 3271
$ws.Range("J132").Value = 4555
$ws.Range("L132").Value = 13665
$ws.Range("N132").Value = -18725

# ----- Sheet: WVR -----
$ws = $wb.Worksheets.Item("WVR")

# Row 107
$ws.Range("H107").Value = 470.4
$ws.Range("J107").Value = 216.5
$ws.Range("L107").Value = 649.5
$ws.Range("N107").Value = -4489.5
# Row 136
$ws.Range("H136").Value = 2180.5
$ws.Range("I136").Value = 1698.2
$ws.Range("K136").Value = 5094.6
$ws.Range("M136").Value = -2544.6
